$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before D (shifts old D:K -> E:L) across the used range
# of the sheet. A bounded range-insert (rather than a whole-column insert) is
# used so the operation doesn't materialise every row of the sheet.
$ws.Range("D5:D102").Insert(-4161)

# Copy the number formats / styles from the column that used to be D (now E)
# into the freshly inserted column D, per financial-statement block, so the
# new cells carry the same date / number styles as the rest of their row.
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)

$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)

$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Populate the new column with the newest fiscal year (FY2018, period ending
# 2018-12-31) figures for the Income Statement, Balance Sheet and Cash Flow
# Statement blocks.

# Income Statement
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 168600
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 37600
$ws.Range("D17").Value = 87100
$ws.Range("D18").Value = 81500
$ws.Range("D20").Value = 73500
$ws.Range("D21").Value = 192600
$ws.Range("D22").Value = "NA"
$ws.Range("D23").Value = 155100
$ws.Range("D24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 155100
$ws.Range("D27").Value = 154400
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -73500
$ws.Range("D33").Value = 154400
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 154400

# Balance Sheet
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 2700
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 109000
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 4000
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 286300
$ws.Range("D48").Value = 1106600
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 5100
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 1513600
$ws.Range("D57").Value = 0
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 35600
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 645000
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 688100
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -37600
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 825500
$ws.Range("D77").Value = 0

# Cash Flow Statement
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 154400
$ws.Range("D83").Value = 37600
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 115500
$ws.Range("D91").Value = -78900
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -3600
$ws.Range("D96").Value = -90400
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -112400
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -400
